$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.694.61'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.283.89'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.81%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.97%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.279.99'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.545'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -10.79%  '
$ws.Range('E11').Value = '  -4.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.508'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.82'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -12.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000248'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.805.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.761.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.273.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.07%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '533.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -8.63%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.114'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -12.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -12.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.761'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -10.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -11.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -10.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -10.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -10.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -11.14%  '
$ws.Range('E31').Value = '  -2.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.09%  '
$ws.Range('E33').Value = '  -15.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -12.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '519.94'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -9.44%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '56.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0446'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0861'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -14.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.127'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -11.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.952.46'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.270'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.83'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -14.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -15.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.114'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '124.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.74%  '
